$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order of writes matters for shared-string table ordering, so we follow
# the same sequence the original author used.

$ws.Range("A13").Value = "Poczatek pasa #2"
$ws.Range("D13").Value = "Koniec pasa #2"

# Update existing labels to reference "pas #1"
$ws.Range("H3").Value = "Roznica latitude dla pas #1"
$ws.Range("H6").Value = "Roznica longitude dla pas #1"

# Update existing formulas to use the new reference row (10) instead of mixed columns
$ws.Range("H4").Formula = "=D4-D10"
$ws.Range("H7").Formula = "=E10-E4"

$ws.Range("H13").Value = "Roznica latitude dla pas #2"

$ws.Range("A14").Value = "Latitude"
$ws.Range("B14").Value = "Longitude"
$ws.Range("D14").Value = "Latitude"
$ws.Range("E14").Value = "Longitude"
$ws.Range("H14").Formula = "=D4-D15"
$ws.Range("I14").Formula = "=H14/5"

$ws.Range("A15").Value = 52.162706
$ws.Range("B15").Value = 20.978023
$ws.Range("D15").Value = 52.170552000000001
$ws.Range("E15").Value = 20.950043999999998

$ws.Range("H16").Value = "Roznica longitude dla pas #2"

$ws.Range("H17").Formula = "=E15-E4"
$ws.Range("I17").Formula = "=H17/5"

$ws.Range("L14").Select()

$wb.Save()
